$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1318
$wsExhibit.Range("F3").Value = 1193
$wsExhibit.Range("F4").Value = 14465
$wsExhibit.Range("F5").Value = 17104
$wsExhibit.Range("F6").Value = 18
$wsExhibit.Range("F7").Value = 140
$wsExhibit.Range("F8").Value = 47
$wsExhibit.Range("F10").Value = 204
$wsExhibit.Range("F16").Value = 41
$wsExhibit.Range("F17").Value = 9
$wsExhibit.Range("F18").Value = 121
$wsExhibit.Range("F19").Value = 0
$wsExhibit.Range("F20").Value = 1302
$wsExhibit.Range("F21").Value = 141
$wsExhibit.Range("F22").Value = 73
$wsExhibit.Range("F25").Value = 7042
$wsExhibit.Range("F27").Value = 30
$wsExhibit.Range("F28").Value = 1153
$wsExhibit.Range("F29").Value = 26
$wsExhibit.Range("F32").Value = 5819
$wsExhibit.Range("F33").Value = 128
$wsExhibit.Range("F35").Value = 214
$wsExhibit.Range("F36").Value = 4976

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1318
$wsAll.Range("F3").Value = 1193
$wsAll.Range("F4").Value = 14465
$wsAll.Range("F5").Value = 17104
$wsAll.Range("F6").Value = 18
$wsAll.Range("F7").Value = 140
$wsAll.Range("F8").Value = 47
$wsAll.Range("F10").Value = 204
$wsAll.Range("F16").Value = 41
$wsAll.Range("F17").Value = 9
$wsAll.Range("F18").Value = 121
$wsAll.Range("F19").Value = 40
$wsAll.Range("F20").Value = 1302
$wsAll.Range("F21").Value = 141
$wsAll.Range("F22").Value = 73
$wsAll.Range("F26").Value = 7042
$wsAll.Range("F28").Value = 30
$wsAll.Range("F29").Value = 1153
$wsAll.Range("F30").Value = 26
$wsAll.Range("F34").Value = 5819
$wsAll.Range("F35").Value = 128
$wsAll.Range("F37").Value = 214
$wsAll.Range("F38").Value = 4976
